$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New week's data is inserted at row 16; the existing rows 16-19 shift down
# to 17-20 (row 19's old data reappears unchanged as the new row 20).
$ws.Rows.Item(16).Insert()

# Row 16 (new record)
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C16").Value = "Ñuble"
$ws.Range("D16").Value = 45161
$ws.Range("D16").NumberFormat = $ws.Range("D17").NumberFormat
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = 100112039
$ws.Range("G16").Value = "Ciboulette"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 2500
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = 2500
$ws.Range("N16").Value = "$/docena de atados"
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 833
$ws.Range("Q16").Value = 3
$ws.Range("R16").Value = "Hortaliza"
